# Sample Project / Main.xlsx - "Rules" sheet
# The rule that used to read "R40" in column B (row 11) is renamed to "1".
# Force the cell to stay text (not be auto-converted to the number 1) by
# switching the cell to Text format before writing the new value - this is
# exactly what happens in Excel when a numeric-looking label is kept as a
# string (the cell keeps its shared-string type, t="s", in the saved file).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
